$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the data (date, quality, volume, prices, unit, etc.) between
# row 4 and row 5, while columns A, B, C, E, F, G, H, O, R stay the same.
# Effectively: row4 gets what used to be in row5, and row5 gets what used to be in row4.

# Read current ("before") values of the columns that change, for row 4 and row 5
$row4_D = $ws.Range("D4").Value()
$row4_I = $ws.Range("I4").Value()
$row4_J = $ws.Range("J4").Value()
$row4_K = $ws.Range("K4").Value()
$row4_L = $ws.Range("L4").Value()
$row4_M = $ws.Range("M4").Value()
$row4_N = $ws.Range("N4").Value()
$row4_P = $ws.Range("P4").Value()
$row4_Q = $ws.Range("Q4").Value()

$row5_D = $ws.Range("D5").Value()
$row5_I = $ws.Range("I5").Value()
$row5_J = $ws.Range("J5").Value()
$row5_K = $ws.Range("K5").Value()
$row5_L = $ws.Range("L5").Value()
$row5_M = $ws.Range("M5").Value()
$row5_N = $ws.Range("N5").Value()
$row5_P = $ws.Range("P5").Value()
$row5_Q = $ws.Range("Q5").Value()

# Write row5's old values into row4
$ws.Range("D4").Value = $row5_D
$ws.Range("I4").Value = $row5_I
$ws.Range("J4").Value = $row5_J
$ws.Range("K4").Value = $row5_K
$ws.Range("L4").Value = $row5_L
$ws.Range("M4").Value = $row5_M
$ws.Range("N4").Value = $row5_N
$ws.Range("P4").Value = $row5_P
$ws.Range("Q4").Value = $row5_Q

# Write row4's old values into row5
$ws.Range("D5").Value = $row4_D
$ws.Range("I5").Value = $row4_I
$ws.Range("J5").Value = $row4_J
$ws.Range("K5").Value = $row4_K
$ws.Range("L5").Value = $row4_L
$ws.Range("M5").Value = $row4_M
$ws.Range("N5").Value = $row4_N
$ws.Range("P5").Value = $row4_P
$ws.Range("Q5").Value = $row4_Q

Write-Output "Swapped row4/row5 values successfully"
